# Apply "All Open as of 20121020" update:
# populate column B (internal venue IDs) for rows 160-196
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    160 = "4eb7a18e82313162702b73b7"
    161 = "4ac518dbf964a52013a920e3"
    162 = "4ea85765e3003a42bcf12fac"
    163 = "4eb0fccc003940a4cd389006"
    164 = "4f47ee93e4b01863527fc7dc"
    165 = "4ed1348430f81894af49af32"
    166 = "4f108084e4b019e990215c3a"
    167 = "4e9d32995503b4a7c128d043"
    168 = "4f594ff3e4b05b5a09964d46"
    169 = "4ea910a7003973c1923a3357"
    170 = "4e9b20d4754aab1ca9dc046c"
    171 = "4e9ef7867beb438c5e78e20b"
    172 = "4e988480d22d4d655ac32ea8"
    173 = "4e98368a02d5a3943708d39d"
    174 = "4e906b718b81c6df65af25bd"
    175 = "4edb9ed930f83fb79e16baa0"
    176 = "4ac518dcf964a52044a920e3"
    177 = "4e8d88ff0039a26265603f19"
    178 = "4e791e988877ce31d210b81e"
    179 = "4e7249e9aeb79ea09553f54c"
    180 = "4e1ed9ddae60cd553deae5cd"
    181 = "4e85ad77e5fab0fb1125fc54"
    182 = "4e653257b993678bb7ca6e65"
    183 = "4d658fad129d60fc8cb814c2"
    185 = "4e569e0b1f6ecd24d016177d"
    186 = "4e3a64ff091a8068d7605b05"
    187 = "4e578dbfb9935d203cf9b155"
    188 = "4e104d56922ec672352e7631"
    189 = "4dfce79d45ddebdfeaac268e"
    190 = "4e53a7e0d4c05e0079f45971"
    193 = "4e14f263d4c0e3d78b28e1e7"
    195 = "4e09e7ab1838eb15af1281c9"
    196 = "4e272ba47d8b5ea7247b1158"
}

$styledRows = @(172, 173, 174, 178, 179, 182, 185, 187)

foreach ($row in $values.Keys | Sort-Object) {
    $cell = $ws.Range("B" + $row)
    $cell.Value = $values[$row]
    if ($styledRows -contains $row) {
        $cell.NumberFormat = "0.00E+00"
    }
}

# Update the view: scroll position and active selection to reflect
# where editing ended (row 179 visible, B196 selected)
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 179
$aw.ScrollColumn = 1
$ws.Range("B196").Select()
